$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 5415.273
$ws.Range("J32").Value = 6788.3335
$ws.Range("L32").Value = 6788.3335
$ws.Range("N32").Value = -7440.3335
# Row 62
$ws.Range("H62").Value = 5376.077
$ws.Range("I62").Value = 4209.8887
$ws.Range("K62").Value = 4209.8887
$ws.Range("M62").Value = -3585.8887
# Row 65
$ws.Range("H65").Value = 5376.077
$ws.Range("I65").Value = 4209.8887
$ws.Range("K65").Value = 21049.4435
$ws.Range("M65").Value = -17929.4435
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 134
$ws.Range("H134").Value = 65939
$ws.Range("J134").Value = 65939
$ws.Range("L134").Value = 65939
$ws.Range("N134").Value = -76079
# Row 137
$ws.Range("H137").Value = 2879.182
$ws.Range("I137").Value = 1549.238
$ws.Range("K137").Value = 4647.714
$ws.Range("M137").Value = -2097.714
# Row 138
$ws.Range("H138").Value = 2636.5974
$ws.Range("I138").Value = 1128.3438
$ws.Range("J138").Value = 3709.1333
$ws.Range("K138").Value = 3385.0314
$ws.Range("L138").Value = 11127.3999
$ws.Range("M138").Value = 1754.9686
$ws.Range("N138").Value = -21407.3999
# Row 140
$ws.Range("H140").Value = 69999
$ws.Range("J140").Value = 69999
$ws.Range("L140").Value = 69999
$ws.Range("N140").Value = -80359

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 2762.25
$ws.Range("I26").Value = 2349.6667
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = 2349.6667
$ws.Range("L26").Value = 4000
$ws.Range("M26").Value = -2019.6667
$ws.Range("N26").Value = -4660
# Row 32
$ws.Range("H32").Value = 2710.0127
$ws.Range("I32").Value = 1961.36
$ws.Range("K32").Value = 1961.36
$ws.Range("M32").Value = -1674.36
# Row 45
$ws.Range("H45").Value = 5198.75
$ws.Range("J45").Value = 2997.5
$ws.Range("L45").Value = 2997.5
$ws.Range("N45").Value = -3751.5
# Row 61
$ws.Range("H61").Value = 5004.2
$ws.Range("I61").Value = 3494.65
$ws.Range("K61").Value = 3494.65
$ws.Range("M61").Value = -3282.65
# Row 136
$ws.Range("H136").Value = 5004.2
$ws.Range("I136").Value = 3494.65
$ws.Range("K136").Value = 10483.95
$ws.Range("M136").Value = -7933.950000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 836.5
$ws.Range("I5").Value = 838
$ws.Range("J5").Value = 835
$ws.Range("K5").Value = 838
$ws.Range("L5").Value = 835
$ws.Range("M5").Value = -725
$ws.Range("N5").Value = -1061
# Row 57
$ws.Range("H57").Value = 54390
$ws.Range("J57").Value = 68780
$ws.Range("L57").Value = 68780
$ws.Range("N57").Value = -70220
# Row 86
$ws.Range("H86").Value = 3742
$ws.Range("I86").Value = 3149.818
$ws.Range("J86").Value = 4556.25
$ws.Range("K86").Value = 3149.818
$ws.Range("L86").Value = 4556.25
$ws.Range("M86").Value = -2026.818
$ws.Range("N86").Value = -6802.25
# Row 89
$ws.Range("H89").Value = 3742
$ws.Range("I89").Value = 3149.818
$ws.Range("J89").Value = 4556.25
$ws.Range("K89").Value = 15749.09
$ws.Range("L89").Value = 22781.25
$ws.Range("M89").Value = -10133.09
$ws.Range("N89").Value = -34013.25
# Row 105
$ws.Range("H105").Value = 13942.131
$ws.Range("I105").Value = 18144.166
$ws.Range("K105").Value = 18144.166
$ws.Range("M105").Value = -16397.166
# Row 136
$ws.Range("H136").Value = 54390
$ws.Range("J136").Value = 68780
$ws.Range("L136").Value = 68780
$ws.Range("N136").Value = -78980

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2407.7778
$ws.Range("I22").Value = 794.6
$ws.Range("J22").Value = 4424.25
$ws.Range("K22").Value = 794.6
$ws.Range("L22").Value = 4424.25
$ws.Range("M22").Value = -444.6
$ws.Range("N22").Value = -5124.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 137.75
$ws.Range("I2").Value = 126.8
$ws.Range("J2").Value = 156
$ws.Range("K2").Value = 760.8
$ws.Range("L2").Value = 936
$ws.Range("M2").Value = -647.8
$ws.Range("N2").Value = -1162
# Row 25
$ws.Range("H25").Value = 169.66667
$ws.Range("J25").Value = 100
$ws.Range("L25").Value = 300
$ws.Range("N25").Value = -638
# Row 30
$ws.Range("H30").Value = 169.66667
$ws.Range("J30").Value = 100
$ws.Range("L30").Value = 300
$ws.Range("N30").Value = -504
# Row 37
$ws.Range("H37").Value = 204444.83
$ws.Range("J37").Value = 204444.83
$ws.Range("L37").Value = 613334.49
$ws.Range("N37").Value = -613558.49
# Row 92
$ws.Range("H92").Value = 4054.9
$ws.Range("I92").Value = 3000
$ws.Range("K92").Value = 9000
$ws.Range("M92").Value = -7752
# Row 136
$ws.Range("H136").Value = 4499.75
$ws.Range("I136").Value = 3999.5
$ws.Range("K136").Value = 11998.5
$ws.Range("M136").Value = -6898.5
# Row 141
$ws.Range("H141").Value = 7157.357
$ws.Range("I141").Value = 3563.3333
$ws.Range("K141").Value = 10689.9999
$ws.Range("M141").Value = -5509.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 668.4375
$ws.Range("I2").Value = 121.2
$ws.Range("J2").Value = 917.1818
$ws.Range("K2").Value = 121.2
$ws.Range("L2").Value = 917.1818
$ws.Range("M2").Value = -8.200000000000003
$ws.Range("N2").Value = -1143.1818

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3232.658
$ws.Range("I46").Value = 1768.1875
$ws.Range("K46").Value = 1768.1875
$ws.Range("M46").Value = -1580.1875
# Row 100
$ws.Range("H100").Value = 10830.048
$ws.Range("J100").Value = 12496.117
$ws.Range("L100").Value = 12496.117
$ws.Range("N100").Value = -13578.117
# Row 122
$ws.Range("H122").Value = 217667.31
$ws.Range("I122").Value = 403226.2
$ws.Range("K122").Value = 1209678.6
$ws.Range("M122").Value = -1207228.6
# Row 132
$ws.Range("H132").Value = 12802.375
$ws.Range("I132").Value = 13851
$ws.Range("J132").Value = 11753.75
$ws.Range("K132").Value = 41553
$ws.Range("L132").Value = 35261.25
$ws.Range("M132").Value = -39023
$ws.Range("N132").Value = -40321.25
# Row 140
$ws.Range("H140").Value = 77429
$ws.Range("J140").Value = 77429
$ws.Range("L140").Value = 77429
$ws.Range("N140").Value = -87789

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 49998.668
$ws.Range("I82").Value = 50000
$ws.Range("J82").Value = 49998
$ws.Range("K82").Value = 50000
$ws.Range("L82").Value = 49998
$ws.Range("M82").Value = -49617
$ws.Range("N82").Value = -50764
# Row 85
$ws.Range("H85").Value = 49998.668
$ws.Range("I85").Value = 50000
$ws.Range("J85").Value = 49998
$ws.Range("K85").Value = 50000
$ws.Range("L85").Value = 49998
$ws.Range("M85").Value = -48674
$ws.Range("N85").Value = -52650
# Row 107
$ws.Range("H107").Value = 1254.3462
$ws.Range("I107").Value = 1401.6818
$ws.Range("K107").Value = 4205.0454
$ws.Range("M107").Value = -2285.0454
# Row 122
$ws.Range("H122").Value = 2264.697
$ws.Range("I122").Value = 1435.875
$ws.Range("J122").Value = 3044.7646
$ws.Range("K122").Value = 4307.625
$ws.Range("L122").Value = 9134.293799999999
$ws.Range("M122").Value = -1857.625
$ws.Range("N122").Value = -14034.2938
